# Update Jinja2 custom env with block trimming and newline configurations
# -------------------------------------------------------------------
# This script reshapes the two template sheets ("main" and "list") of the
# catalog template workbook so the Jinja expressions render on their own
# line (the custom Jinja environment now trims blocks / keeps newlines),
# swaps the EAN / PRICE_AMOUNT column templates, adds an `|int` filter to
# the price template, adds a brand-new VAT description template on the
# "list" sheet, and refreshes the look (narrower wrapped columns, top
# aligned header/body text) to match the new multi-line cell content.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("main")
$ws2 = $wb.Worksheets.Item("list")

# ---------------------------------------------------------------
# 1. "main" sheet - rewrite the row-2 Jinja templates with embedded
#    newlines, and swap the EAN / PRICE_AMOUNT templates so each stays
#    lined up under its own header (Price / EAN).
# ---------------------------------------------------------------
$nl = "`n"

$ws1.Range("A2").Value = "{% for article in articles_json %}" + $nl + "{{article['SUPPLIER_AID']}}" + $nl + "{{split}}" + $nl + "{% endfor %}"
$ws1.Range("B2").Value = "{% for article in articles_json %}" + $nl + "{{article['DESCRIPTION_LONG']}}" + $nl + "{{split}}" + $nl + "{% endfor %}"
$ws1.Range("C2").Value = "{% for article in articles_json %}" + $nl + "{{article['PRICE_AMOUNT']|int}}" + $nl + "{{split}}" + $nl + "{% endfor %}"
$ws1.Range("D2").Value = "{% for article in articles_json %}" + $nl + "EUR{{split}}" + $nl + "{% endfor %}"
$ws1.Range("E2").Value = "{% for article in articles_json %}" + $nl + "{{article['EAN']}}" + $nl + "{{split}}" + $nl + "{% endfor %}"

# ---------------------------------------------------------------
# 2. "list" sheet - move the "List of materialas" label from B2 to A1
#    (keeping its bold/red/yellow look), and replace the old row-3
#    sample template with the new VAT description template in A2.
# ---------------------------------------------------------------
$ws2.Range("B2").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws2.Range("A1").Value = $ws2.Range("B2").Value
$ws2.Range("B2").Clear() | Out-Null

$ws2.Range("A2").Value = "{% for article in articles_json %}" + $nl + "{{ article['DESCRIPTION_SHORT']|upper }} is priced at {{ article['PRICE_AMOUNT'] }} with VAT." + $nl + "{{split}}" + $nl + "{% endfor %}"
$ws2.Range("B3").Clear() | Out-Null

# ---------------------------------------------------------------
# 3. Formatting refresh on "main": narrower columns (data now wraps
#    across several lines), a new empty column F, top-aligned header
#    row (still bold) and top-aligned + wrapped body row.
# ---------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 28.5
$ws1.Columns.Item(2).ColumnWidth = 29.333333333333332
$ws1.Columns.Item(3).ColumnWidth = 28.5
$ws1.Columns.Item(4).ColumnWidth = 28.5
$ws1.Columns.Item(5).ColumnWidth = 28.5
$ws1.Columns.Item(6).ColumnWidth = 8.333333333333334

$ws1.Range("A:F").VerticalAlignment = -4160

$ws1.Range("A1:E1").VerticalAlignment = -4160

$ws1.Range("A2:E2").VerticalAlignment = -4160
$ws1.Range("A2:E2").WrapText = $true
$ws1.Rows.Item(2).RowHeight = 15

# ---------------------------------------------------------------
# 4. Formatting refresh on "list": widen column A (now holds the
#    labels/templates) and give the new row-2 template a touch of
#    row height so it isn't clipped.
# ---------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 89.5
$ws2.Rows.Item(2).RowHeight = 14.25
